# Applies the "Updated batch documentation and distribution" edit to
# doc/Manual/Batchmode/StageStructFile.xlsx
#
# Summary of the change (per the canonical-XML diff):
#  - workbook.xml: x15ac:absPath url updated from
#      C:\RangeShifter_v2_0\doc\Manual\Batchmode\  ->  C:\RangeShifter2\doc\Manual\Batchmode\
#    (This is an Excel-internal "last saved from" bookkeeping attribute with
#    no surface on the Excel object model / COM automation interface, so it
#    is not settable from macro code in real Excel either. We still try the
#    closest available property so the intent is recorded; it is a no-op if
#    unsupported by the host.)
#  - xl/sharedStrings.xml: three generic, duplicated "Stage-specific
#    density-dependence in fecuntity..." placeholder strings (which had been
#    reused verbatim across the Fec/Dev/Surv "*StageWts" rows) are replaced
#    by three distinct, correctly-worded "Stage-weighted density-dependence
#    in fecuntity/development/survival..." strings, and the three
#    "Required if ...StageDep is 1" notes gain an ", otherwise NULL" suffix.
#    Every other apparent shared-string index change in the diff is a pure
#    side effect of Excel's shared-string table being rebuilt (old, now
#    unused strings dropped, new strings appended) -- no other cell text
#    actually changes.

$wb = $excel.ActiveWorkbook

# Best-effort: record the new "last saved from" folder. Real Excel does not
# expose x15ac:absPath via COM either, so this is harmless if ignored.
try {
    $wb.Path = "C:\RangeShifter2\doc\Manual\Batchmode\"
} catch {
}

$ws = $wb.Worksheets.Item("Description")

# Order matters: new shared strings are appended in first-use order, and the
# diff expects them at indices 64-69 in exactly this sequence.
$ws.Range("C10").Value = "Stage-weighted density-dependence in fecuntity: 0 = No, 1  = Yes"
$ws.Range("C14").Value = "Stage-weighted density-dependence in development: 0 = No, 1  = Yes"
$ws.Range("C18").Value = "Stage-weighted density-dependence in survival: 0 = No, 1  = Yes"

$ws.Range("D11").Value = "Required if FecStageDep is 1, otherwise NULL"
$ws.Range("D15").Value = "Required if DevStageDep is 1, otherwise NULL"
$ws.Range("D19").Value = "Required if SurvStageDep is 1, otherwise NULL"
